$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999469869549,
    0.99487436656238937,
    0.97712642219323398,
    0.96988259040409908,
    0.95962967880609273,
    0.93752843719426737,
    0.93736854694320182,
    0.93345579283975955,
    0.93284193401942517,
    0.9331717881118029,
    0.9328682403146763,
    0.93254676384868884,
    0.93192491333540461,
    0.93212471716385159,
    0.9338974299706635,
    0.93671529255934805,
    0.94343993763906875,
    0.94763497246407669,
    0.99142518447265959,
    0.98430804352627843,
    0.98290954510062678,
    0.98164503659659441,
    0.97009108286373424,
    0.95706965653353326,
    0.9506125679711932,
    0.93387237340636398,
    0.92902509633214048,
    0.90754985977928848,
    0.89227755637408579,
    0.88570694052820653,
    0.87805292811573887,
    0.87637358095193219,
    0.87585356764055822
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
